$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: cpsc 261 204 L2A " " " " WED 1200 1300
$ws.Range("A4").Value = "cpsc"
$ws.Range("B4").Value = 261
$ws.Range("C4").Value = 204
$ws.Range("D4").Value = "L2A"
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = "WED"
$ws.Range("H4").Value = 1200
$ws.Range("I4").Value = 1300

# Row 5: comm 280 201 " " " " D2B FRI 2000 2100
$ws.Range("A5").Value = "comm"
$ws.Range("B5").Value = 280
$ws.Range("C5").Value = 201
$ws.Range("D5").Value = " "
$ws.Range("E5").Value = " "
$ws.Range("F5").Value = "D2B"
$ws.Range("G5").Value = "FRI"
$ws.Range("H5").Value = 2000
$ws.Range("I5").Value = 2100

# Update the active selection to match the post-edit state
$ws.Range("K6").Select()
